$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"/"heading"/"Outstanding"
# shift one column to the right), copying the width of the column to its
# left (M, "In Advance") the way Excel's Insert Column does.
$mWidth = $ws.Columns("M").ColumnWidth
$null = $ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with K13 selected.
$ws.Activate()
$null = $ws.Range("K13").Select()
